# TD-6649 add Business Concept file manager domain name field
#
# The upload template's second header column ("domain") is being repurposed
# as "domain_external_id" so a separate "domain" (name) field can be added
# later; for this fixture update only the header cell text changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "domain_external_id"

# Reflect that the user last interacted with the renamed header cell.
$ws.Range("B1").Select() | Out-Null
